$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 366.5
$ws.Range("I12").Value = 322
$ws.Range("J12").Value = 500
$ws.Range("K12").Value = 322
$ws.Range("L12").Value = 500
$ws.Range("M12").Value = -152
$ws.Range("N12").Value = -840
$ws.Range("H40").Value = 1000
$ws.Range("J40").Value = 1000
$ws.Range("L40").Value = 1000
$ws.Range("N40").Value = -1350
$ws.Range("H111").Value = 1200
$ws.Range("I111").Value = 1200
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 3600
$ws.Range("L111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -533
$ws.Range("H121").Value = 2055.8096
$ws.Range("J121").Value = 2055.8096
$ws.Range("L121").Value = 6167.4288
$ws.Range("N121").Value = -9661.4288
$ws.Range("H125").Value = 1676
$ws.Range("I125").Value = 866
$ws.Range("K125").Value = 7794
$ws.Range("M125").Value = -5334
$ws.Range("H129").Value = 2267
$ws.Range("J129").Value = 1499.5
$ws.Range("L129").Value = 4498.5
$ws.Range("N129").Value = -14498.5
$ws.Range("H132").Value = 1500.3334
$ws.Range("I132").Value = 1500.3334
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4501.0002
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -1971.0002
$ws.Range("H135").Value = 2434.1667
$ws.Range("I135").Value = 2434.1667
$ws.Range("K135").Value = 21907.5003
$ws.Range("M135").Value = -19372.5003
$ws.Range("H137").Value = 1687.2452
$ws.Range("I137").Value = 1461.9048
$ws.Range("K137").Value = 4385.7144
$ws.Range("M137").Value = -1835.7144

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15196.025
$ws.Range("I32").Value = 14806.595
$ws.Range("K32").Value = 14806.595
$ws.Range("M32").Value = -14519.595
$ws.Range("H74").Value = 1582.2424
$ws.Range("I74").Value = 1084.72
$ws.Range("K74").Value = 1084.72
$ws.Range("M74").Value = -210.72
$ws.Range("H77").Value = 1582.2424
$ws.Range("I77").Value = 1084.72
$ws.Range("K77").Value = 5423.6
$ws.Range("M77").Value = -1055.6
$ws.Range("H132").Value = 2702.5833
$ws.Range("I132").Value = 2094
$ws.Range("J132").Value = 3554.6
$ws.Range("K132").Value = 6282
$ws.Range("L132").Value = 10663.8
$ws.Range("M132").Value = -3752
$ws.Range("N132").Value = -15723.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2582.16
$ws.Range("I86").Value = 2533.8
$ws.Range("J86").Value = 2775.6
$ws.Range("K86").Value = 2533.8
$ws.Range("L86").Value = 2775.6
$ws.Range("M86").Value = -1410.8
$ws.Range("N86").Value = -5021.6
$ws.Range("H89").Value = 2582.16
$ws.Range("I89").Value = 2533.8
$ws.Range("J89").Value = 2775.6
$ws.Range("K89").Value = 12669
$ws.Range("L89").Value = 13878
$ws.Range("M89").Value = -7053
$ws.Range("N89").Value = -25110

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2225.2173
$ws.Range("I31").Value = 2204.6667
$ws.Range("J31").Value = 2441
$ws.Range("K31").Value = 2204.6667
$ws.Range("L31").Value = 2441
$ws.Range("M31").Value = -1909.6667
$ws.Range("N31").Value = -3031
$ws.Range("H34").Value = 2225.2173
$ws.Range("I34").Value = 2204.6667
$ws.Range("J34").Value = 2441
$ws.Range("K34").Value = 2204.6667
$ws.Range("L34").Value = 2441
$ws.Range("M34").Value = -2002.6667
$ws.Range("N34").Value = -2845
$ws.Range("H50").Value = 50092
$ws.Range("J50").Value = 50092
$ws.Range("L50").Value = 50092
$ws.Range("N50").Value = -51342
$ws.Range("H51").Value = 50099
$ws.Range("J51").Value = 50099
$ws.Range("L51").Value = 50099
$ws.Range("N51").Value = -51571
$ws.Range("H59").Value = 60127
$ws.Range("J59").Value = 60127
$ws.Range("L59").Value = 60127
$ws.Range("N59").Value = -62417
$ws.Range("H60").Value = 36381.87
$ws.Range("J60").Value = 36853.773
$ws.Range("L60").Value = 36853.773
$ws.Range("N60").Value = -37875.773
$ws.Range("H61").Value = 50099
$ws.Range("J61").Value = 50099
$ws.Range("L61").Value = 50099
$ws.Range("N61").Value = -50795
$ws.Range("H86").Value = 22534.8
$ws.Range("J86").Value = 34474.5
$ws.Range("L86").Value = 34474.5
$ws.Range("N86").Value = -36720.5
$ws.Range("H89").Value = 22534.8
$ws.Range("J89").Value = 34474.5
$ws.Range("L89").Value = 172372.5
$ws.Range("N89").Value = -183604.5
$ws.Range("H99").Value = 9334.571
$ws.Range("I99").Value = 9085.75
$ws.Range("K99").Value = 9085.75
$ws.Range("M99").Value = -7587.75
$ws.Range("H126").Value = 9334.571
$ws.Range("I126").Value = 9085.75
$ws.Range("K126").Value = 27257.25
$ws.Range("M126").Value = -24787.25
$ws.Range("H132").Value = 4698.65
$ws.Range("I132").Value = 4665.278
$ws.Range("K132").Value = 13995.834
$ws.Range("M132").Value = -11465.834
$ws.Range("H134").Value = 2399.5881
$ws.Range("J134").Value = 4998
$ws.Range("L134").Value = 14994
$ws.Range("N134").Value = -20064

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 216.4
$ws.Range("I7").Value = 195.5
$ws.Range("J7").Value = 300
$ws.Range("K7").Value = 586.5
$ws.Range("L7").Value = 900
$ws.Range("M7").Value = -474.5
$ws.Range("N7").Value = -1124
$ws.Range("H33").Value = 392.5
$ws.Range("I33").Value = 469
$ws.Range("K33").Value = 2814
$ws.Range("M33").Value = -2531
$ws.Range("H122").Value = 112977.445
$ws.Range("J122").Value = 144828.14
$ws.Range("L122").Value = 1303453.26
$ws.Range("N122").Value = -1308353.26
$ws.Range("H132").Value = 3862.2856
$ws.Range("I132").Value = 2439.6
$ws.Range("J132").Value = 4652.6665
$ws.Range("K132").Value = 21956.4
$ws.Range("L132").Value = 41873.9985
$ws.Range("M132").Value = -19426.4
$ws.Range("N132").Value = -46933.9985

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 697.76
$ws.Range("I97").Value = 620.6667
$ws.Range("J97").Value = 1102.5
$ws.Range("K97").Value = 620.6667
$ws.Range("L97").Value = 1102.5
$ws.Range("M97").Value = -124.6667
$ws.Range("N97").Value = -2094.5
$ws.Range("H132").Value = 4998.75
$ws.Range("I132").Value = 4998
$ws.Range("K132").Value = 14994
$ws.Range("M132").Value = -12464

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1450.75
$ws.Range("I7").Value = 1450.75
$ws.Range("K7").Value = 1450.75
$ws.Range("M7").Value = -1338.75
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H40").Value = 2998
$ws.Range("I40").Value = 2998
$ws.Range("K40").Value = 2998
$ws.Range("M40").Value = -2862
$ws.Range("H46").Value = 449
$ws.Range("I46").Value = 449
$ws.Range("K46").Value = 449
$ws.Range("M46").Value = -261
$ws.Range("H61").Value = 3826.375
$ws.Range("I61").Value = 3730.1428
$ws.Range("K61").Value = 3730.1428
$ws.Range("M61").Value = -3528.1428
$ws.Range("H100").Value = 1907.8889
$ws.Range("I100").Value = 1224.4
$ws.Range("K100").Value = 1224.4
$ws.Range("M100").Value = -683.4000000000001
$ws.Range("H113").Value = 3826.375
$ws.Range("I113").Value = 3730.1428
$ws.Range("K113").Value = 3730.1428
$ws.Range("M113").Value = -1560.1428
$ws.Range("H126").Value = 1450.75
$ws.Range("I126").Value = 1450.75
$ws.Range("K126").Value = 4352.25
$ws.Range("M126").Value = -1882.25
$ws.Range("H132").Value = 5874.4
$ws.Range("I132").Value = 5791.3335
$ws.Range("K132").Value = 17374.0005
$ws.Range("M132").Value = -14844.0005
$ws.Range("H136").Value = 6498.0625
$ws.Range("I136").Value = 6475.231
$ws.Range("K136").Value = 19425.693
$ws.Range("M136").Value = -16875.693

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 12122.125
$ws.Range("I81").Value = 6994
$ws.Range("J81").Value = 17250.25
$ws.Range("K81").Value = 13988
$ws.Range("L81").Value = 34500.5
$ws.Range("M81").Value = -12927
$ws.Range("N81").Value = -36622.5
$ws.Range("H84").Value = 12122.125
$ws.Range("I84").Value = 6994
$ws.Range("J84").Value = 17250.25
$ws.Range("K84").Value = 69940
$ws.Range("L84").Value = 172502.5
$ws.Range("M84").Value = -64636
$ws.Range("N84").Value = -183110.5
$ws.Range("H96").Value = 986.3333
$ws.Range("I96").Value = 986.3333
$ws.Range("K96").Value = 986.3333
$ws.Range("M96").Value = 386.6667
$ws.Range("H100").Value = 1297.1666
$ws.Range("I100").Value = 1297.1666
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2594.3332
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -2053.3332
$ws.Range("H107").Value = 603
$ws.Range("I107").Value = 546.7143
$ws.Range("K107").Value = 1640.1429
$ws.Range("M107").Value = 279.8571000000002
$ws.Range("H136").Value = 2185.1428
$ws.Range("I136").Value = 1766.3334
$ws.Range("J136").Value = 2499.25
$ws.Range("K136").Value = 5299.0002
$ws.Range("L136").Value = 7497.75
$ws.Range("M136").Value = -2749.0002
$ws.Range("N136").Value = -12597.75
